$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 349.5
$ws.Range("I61").Value = 349.5
$ws.Range("K61").Value = 1048.5
$ws.Range("M61").Value = -876.5
$ws.Range("H98").Value = 1479.25
$ws.Range("I98").Value = 1479.25
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1479.25
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 18.75
$ws.Range("N98").ClearContents()
$ws.Range("H122").Value = 1479.25
$ws.Range("I122").Value = 1479.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4437.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1987.75
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 4955
$ws.Range("I132").Value = 5348.8887
$ws.Range("J132").Value = 3942.1428
$ws.Range("K132").Value = 16046.6661
$ws.Range("L132").Value = 11826.4284
$ws.Range("M132").Value = -13516.6661
$ws.Range("N132").Value = -16886.4284
$ws.Range("H138").Value = 1512.7142
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 1512.7142
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 4538.142599999999
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -14818.1426
$ws.Range("H141").Value = 12338.8
$ws.Range("I141").Value = 12338.8
$ws.Range("K141").Value = 37016.39999999999
$ws.Range("M141").Value = -31836.39999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6644.1934
$ws.Range("I32").Value = 5895.552
$ws.Range("K32").Value = 5895.552
$ws.Range("M32").Value = -5608.552
$ws.Range("H45").Value = 2299.2
$ws.Range("I45").Value = 2124.25
$ws.Range("K45").Value = 2124.25
$ws.Range("M45").Value = -1747.25
$ws.Range("H88").Value = 2970
$ws.Range("I88").Value = 2909
$ws.Range("J88").Value = 3010.6667
$ws.Range("K88").Value = 2909
$ws.Range("L88").Value = 3010.6667
$ws.Range("M88").Value = -2503
$ws.Range("N88").Value = -3822.6667
$ws.Range("H91").Value = 2970
$ws.Range("I91").Value = 2909
$ws.Range("J91").Value = 3010.6667
$ws.Range("K91").Value = 2909
$ws.Range("L91").Value = 3010.6667
$ws.Range("M91").Value = -1505
$ws.Range("N91").Value = -5818.6667
$ws.Range("H102").Value = 4505.5
$ws.Range("I102").Value = 5000
$ws.Range("K102").Value = 5000
$ws.Range("M102").Value = -3378
$ws.Range("H122").Value = 1580.7333
$ws.Range("I122").Value = 1481.5
$ws.Range("J122").Value = 1977.6666
$ws.Range("K122").Value = 4444.5
$ws.Range("L122").Value = 5932.9998
$ws.Range("M122").Value = -1994.5
$ws.Range("N122").Value = -10832.9998
$ws.Range("H139").Value = 99999.336
$ws.Range("J139").Value = 99999.336
$ws.Range("L139").Value = 99999.336
$ws.Range("N139").Value = -110279.336

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2808.5
$ws.Range("I105").Value = 2731.3333
$ws.Range("K105").Value = 2731.3333
$ws.Range("M105").Value = -984.3332999999998

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1970.3334
$ws.Range("I58").Value = 1220
$ws.Range("J58").Value = 2506.2856
$ws.Range("K58").Value = 1220
$ws.Range("L58").Value = 2506.2856
$ws.Range("M58").Value = -1017
$ws.Range("N58").Value = -2912.2856
$ws.Range("H136").Value = 1970.3334
$ws.Range("I136").Value = 1220
$ws.Range("J136").Value = 2506.2856
$ws.Range("K136").Value = 3660
$ws.Range("L136").Value = 7518.8568
$ws.Range("M136").Value = -1110
$ws.Range("N136").Value = -12618.8568

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 600
$ws.Range("I18").Value = 600
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 1800
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -1631
$ws.Range("N18").ClearContents()
$ws.Range("H68").Value = 1100
$ws.Range("I68").Value = 700
$ws.Range("J68").Value = 1500
$ws.Range("K68").Value = 2100
$ws.Range("L68").Value = 4500
$ws.Range("M68").Value = -1289
$ws.Range("N68").Value = -6122
$ws.Range("H71").Value = 1100
$ws.Range("I71").Value = 700
$ws.Range("J71").Value = 1500
$ws.Range("K71").Value = 6300
$ws.Range("L71").Value = 13500
$ws.Range("M71").Value = -2244
$ws.Range("N71").Value = -21612
$ws.Range("H81").Value = 1499
$ws.Range("I81").Value = 1499
$ws.Range("K81").Value = 4497
$ws.Range("M81").Value = -3374
$ws.Range("H84").Value = 1499
$ws.Range("I84").Value = 1499
$ws.Range("K84").Value = 13491
$ws.Range("M84").Value = -7875
$ws.Range("H113").Value = 2522
$ws.Range("J113").Value = 2533
$ws.Range("L113").Value = 7599
$ws.Range("N113").Value = -11939
$ws.Range("H117").Value = 2999
$ws.Range("I117").Value = 2999
$ws.Range("K117").Value = 8997
$ws.Range("M117").Value = -5555
$ws.Range("H119").Value = 1687
$ws.Range("I119").Value = 1687
$ws.Range("K119").Value = 5061
$ws.Range("M119").Value = -223
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H133").Value = 18948.75
$ws.Range("I133").Value = 17931.666
$ws.Range("K133").Value = 53794.99800000001
$ws.Range("M133").Value = -48734.99800000001
$ws.Range("H134").Value = 5515
$ws.Range("I134").Value = 5515
$ws.Range("K134").Value = 16545
$ws.Range("M134").Value = -11475
$ws.Range("H139").Value = 9795
$ws.Range("I139").Value = 9795
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 29385
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -24245
$ws.Range("N139").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 2695949.8
$ws.Range("I3").Value = 5003233
$ws.Range("K3").Value = 5003233
$ws.Range("M3").Value = -5003117
$ws.Range("H80").Value = 4999.75
$ws.Range("I80").Value = 4999
$ws.Range("K80").Value = 4999
$ws.Range("M80").Value = -4001
$ws.Range("H83").Value = 4999.75
$ws.Range("I83").Value = 4999
$ws.Range("K83").Value = 24995
$ws.Range("M83").Value = -20003
$ws.Range("H102").Value = 1751.3182
$ws.Range("I102").Value = 1751.3182
$ws.Range("K102").Value = 1751.3182
$ws.Range("M102").Value = -129.3181999999999
$ws.Range("H107").Value = 2800
$ws.Range("I107").Value = 400.33334
$ws.Range("K107").Value = 400.33334
$ws.Range("M107").Value = 1519.66666
$ws.Range("H122").Value = 2558
$ws.Range("J122").Value = 3336.3333
$ws.Range("L122").Value = 10008.9999
$ws.Range("N122").Value = -14908.9999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3833.3333
$ws.Range("I22").Value = 750
$ws.Range("K22").Value = 750
$ws.Range("M22").Value = -455
$ws.Range("H27").Value = 3833.3333
$ws.Range("I27").Value = 750
$ws.Range("K27").Value = 750
$ws.Range("M27").Value = -643
$ws.Range("H40").Value = 3573.25
$ws.Range("J40").Value = 4999.5
$ws.Range("L40").Value = 4999.5
$ws.Range("N40").Value = -5271.5
$ws.Range("H93").Value = 3404.2856
$ws.Range("I93").Value = 3446.5
$ws.Range("K93").Value = 3446.5
$ws.Range("M93").Value = -2198.5
$ws.Range("H122").Value = 3749.75
$ws.Range("I122").Value = 3666.6667
$ws.Range("J122").Value = 3999
$ws.Range("K122").Value = 11000.0001
$ws.Range("L122").Value = 11997
$ws.Range("M122").Value = -8550.000100000001
$ws.Range("N122").Value = -16897
$ws.Range("H132").Value = 3501.6155
$ws.Range("I132").Value = 3554.8948
$ws.Range("J132").Value = 3357
$ws.Range("K132").Value = 10664.6844
$ws.Range("L132").Value = 10071
$ws.Range("M132").Value = -8134.6844
$ws.Range("N132").Value = -15131
$ws.Range("H136").Value = 6467.8
$ws.Range("I136").Value = 4099.143
$ws.Range("K136").Value = 12297.429
$ws.Range("M136").Value = -9747.429

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 574.3333
$ws.Range("I113").Value = 361.5
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1084.5
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1085.5
$ws.Range("N113").Value = -7340
$ws.Range("H122").Value = 3713.3572
$ws.Range("I122").Value = 3589.2727
$ws.Range("J122").Value = 4168.3335
$ws.Range("K122").Value = 10767.8181
$ws.Range("L122").Value = 12505.0005
$ws.Range("M122").Value = -8317.8181
$ws.Range("N122").Value = -17405.0005
$ws.Range("H136").Value = 27840.2
$ws.Range("I136").Value = 16672.428
$ws.Range("J136").Value = 53898.332
$ws.Range("K136").Value = 50017.284
$ws.Range("L136").Value = 161694.996
$ws.Range("M136").Value = -47467.284
$ws.Range("N136").Value = -166794.996
